$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "C2"
$ws.Range("E1").Value = "C2*2"

$ws.Range("D2").Formula = "=C2"
$ws.Range("E2").Formula = "=C2*2"

$ws.Range("E3").Select()
